$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab_7a_Quellen")

# Fill column E ("LinkMeldungEn") with source attribution text for each row
$ws.Range("E2").Value = "the Federal Foreign Office "
$ws.Range("E3").Value = "the AG Energiebilanzen e.V."
$ws.Range("E4").Value = "the Working Group on Renewable Energy Statistics"
$ws.Range("E5").Value = "the Sustainable Development Working Group"
$ws.Range("E6").Value = "the Arbeitskreis der Umweltökonomischen Gesamtrechnungen der Länder"
$ws.Range("E7").Value = "the Agricultural Market Information Company"
$ws.Range("E8").Value = "the Federal Institute for Research on Building, Urban Affairs and Spatial Development"
$ws.Range("E9").Value = "the Federal Agency for Nature Conservation "
$ws.Range("E10").Value = "the Federal Criminal Police Office"
$ws.Range("E11").Value = "the Federal Chancellery"
$ws.Range("E12").Value = "the Federal Government Commissioner for Culture and the Media"
$ws.Range("E13").Value = "the Federal Ministry of Education and Research"
$ws.Range("E14").Value = "the Federal Ministry for Digital and Transport"
$ws.Range("E15").Value = "the Federal Ministry of Food and Agriculture"
$ws.Range("E16").Value = "the Federal Ministry of Finance"
$ws.Range("E17").Value = "the Federal Ministry of Health"
$ws.Range("E18").Value = "the Federal Ministry for Economic Affairs and Climate Action"
$ws.Range("E19").Value = "the Federal Ministry for Economic Cooperation and Development"
$ws.Range("E20").Value = "the Organic Food Production Alliance"
$ws.Range("E21").Value = "the German Federal Bank"
$ws.Range("E22").Value = "the Federal Centre for Health Education"
$ws.Range("E23").Value = "the CEval GmbH"
$ws.Range("E25").Value = "the German EMAS Advisory Board"
$ws.Range("E26").Value = "the German Institute for Economic Research"
$ws.Range("E27").Value = "the European Commission"
$ws.Range("E28").Value = "Eurostat"
$ws.Range("E29").Value = "the European Central Bank"
$ws.Range("E30").Value = "the Frauen in die Aufsichtsräte e.V."
$ws.Range("E31").Value = "the Federal Health Monitoring"
$ws.Range("E32").Value = "the Gesellschaft für Konsumforschung"
$ws.Range("E33").Value = "the Deutsche Gesellschaft für Internationale Zusammenarbeit GmbH"
$ws.Range("E34").Value = "the Institute for Energy and Environmental Research"
$ws.Range("E35").Value = "the Johann Heinrich von Thünen Institute"
$ws.Range("E36").Value = "the Institute for Crop and Soil Science, Julius Kühn Institute"
$ws.Range("E37").Value = "the Federal Motor Transport Authority"
$ws.Range("E38").Value = "the Kreditanstalt für Wiederaufbau"
$ws.Range("E39").Value = "the Center of Excellence for Sustainable Procurement"
$ws.Range("E41").Value = "the Länderinitiative Kernindikatoren"
$ws.Range("E42").Value = "the Organisation for Economic Co-operation and Development"
$ws.Range("E43").Value = "the Robert Koch Institute"
$ws.Range("E44").Value = "the Stifterverband Wissenschaftsstatistik"
$ws.Range("E45").Value = "Transparency International"
$ws.Range("E46").Value = "the German Environment Agency"
$ws.Range("E47").Value = "the German Environment Agency"
$ws.Range("E48").Value = "the German Environment Agency"
$ws.Range("E49").Value = "the University of Giessen"
$ws.Range("E50").Value = "the Verkehrsclub Deutschland e.V."
$ws.Range("E51").Value = "the Statistische Ämter des Bundes und der Länder"
$ws.Range("E52").Value = "the World Health Organization"

# Fix D47: add missing closing parenthesis
$ws.Range("D47").Value = "German Environment Agency (as reported by the Länder and by river basin commissions)"

$wb.Save()
